$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows (13 and 14) below the existing data, copying the
# formatting (styles) from row 12 so the new blank rows look the same
# as the others (s="2" on col A, s="1" on cols B/C), but without values.
$ws.Range("A12:C12").Copy()
$ws.Range("A13:C14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Clear the contents (values) of every row except row 5, keeping the
# cell formatting/styles intact.
$ws.Range("A1:C4").ClearContents()
$ws.Range("A6:C14").ClearContents()

# Row 5 is the only row that keeps data; update its quantity to 0.
$ws.Range("C5").Value = 0

# Update the active selection to reflect the last-edited cell.
$ws.Range("C5").Select()
